$d = $word.ActiveDocument

# Target the specific sentence about the 90% intensity check (there is another,
# unrelated "...samples)..." elsewhere in the document, so match enough context
# to be unambiguous).
$old = "in at least 90% of samples). This intensity check minimizes confusion"
$new = "in at least 90% of samples, including pooled samples). This intensity check minimizes confusion"

$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

if (-not $found) {
    throw "Target text for the pooled-samples edit was not found."
}

Write-Host "Replaced: $found"
